# Session-27 update: remove the "Google Maps Demo" slide (slide 2).
# All other slides shift up by one position; PowerPoint's own save
# logic takes care of renumbering relationship ids, notes-slide
# field caches, and similar bookkeeping automatically.

$p = $ppt.ActivePresentation
$p.Slides.Item(2).Delete()
